# Applies the PS-VRP "Problemi_lettura_Excel" data refresh described in the
# commit "Risoluzione bug check nelle ricerche locali".
#
# The sheet holds a 10-row extract (header + 9 data rows... actually header
# row 1 + rows 2-10). After the refresh it only has 7 data rows (rows 2-8),
# most values shifted/changed, and the used range shrinks to A1:P8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the last two data rows (10 -> 8 total rows) -------------------
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()

# --- 2. Overwrite the changed cell values in rows 2-8 -----------------------
$ws.Range("B2").Value = 45903
$ws.Range("A3").Value = 251702
$ws.Range("B3").Value = 45903
$ws.Range("D3").Value = 45853
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 60063
$ws.Range("G3").Value = 4607
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1"
$ws.Range("H3").ClearFormats()
$ws.Range("J3").Value = 340
$ws.Range("K3").Value = 1020
$ws.Range("A4").Value = 251310
$ws.Range("B4").Value = 45770
$ws.Range("C4").Value = "STAMPATO"
$ws.Range("D4").Value = 45769
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 8611
$ws.Range("G4").Value = 336
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "6"
$ws.Range("H4").ClearFormats()
$ws.Range("J4").Value = 410
$ws.Range("K4").Value = 820
$ws.Range("L4").Value = "CAMPO VUOTO"
$ws.Range("N4").Value = 1
$ws.Range("A5").Value = 252683
$ws.Range("B5").Value = 45903
$ws.Range("C5").Value = "STAMPATO"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 19381
$ws.Range("G5").Value = 2000
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "5"
$ws.Range("H5").ClearFormats()
$ws.Range("J5").Value = 820
$ws.Range("K5").Value = 820
$ws.Range("L5").Value = "CAMPO VUOTO"
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = "X"
$ws.Range("P5").Value = 40295
$ws.Range("A6").Value = 252459
$ws.Range("B6").Value = 45855
$ws.Range("D6").Value = "CAMPO VUOTO"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 4000
$ws.Range("G6").Value = 91
$ws.Range("H6").Value = "CAMPO VUOTO"
$ws.Range("I6").Value = "foglio"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 480
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = "CAMPO VUOTO"
$ws.Range("P6").Value = "CAMPO VUOTO"
$ws.Range("A7").Value = 252596
$ws.Range("B7").Value = 45904
$ws.Range("D7").Value = 45873
$ws.Range("F7").Value = 3497
$ws.Range("G7").Value = 305
$ws.Range("H7").Value = "5 / 9"
$ws.Range("I7").Value = "bobina"
$ws.Range("J7").Value = 470
$ws.Range("K7").Value = 940
$ws.Range("M7").Value = "Dati OK"
$ws.Range("O7").Value = "X"
$ws.Range("P7").Value = 40295
$ws.Range("A8").Value = 251624
$ws.Range("B8").Value = 45908
$ws.Range("D8").Value = 45895
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 10957
$ws.Range("G8").Value = 574
$ws.Range("H8").Value = "CAMPO VUOTO"
$ws.Range("J8").Value = 155
$ws.Range("K8").Value = 620
$ws.Range("N8").Value = 0
$ws.Range("P8").Value = 40308

# --- 3. Fix up formatting for cells whose value "kind" changed -------------
# B4 / B5 become real dates again (were showing as "CAMPO VUOTO" placeholder
# text) - reuse the existing date-formatted style from B2.
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)

# L4 / L5 become "CAMPO VUOTO" placeholders - reuse the existing highlighted
# style from D2.
$ws.Range("D2").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L5").PasteSpecial(-4122)

# D6 becomes a "CAMPO VUOTO" placeholder (was a formatted date).
$ws.Range("D6").PasteSpecial(-4122)

# L6 and I7 go back to plain, unformatted cells.
$ws.Range("L6").ClearFormats()
$ws.Range("I7").ClearFormats()

$excel.CutCopyMode = 0
